# Testing y arreglos. Front.
# Adds new test rows (53-83) to the "Testing y funcionalidades" sheet,
# widens several columns, and moves the active selection to reflect the
# extended data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting (fill colors) for new rows from existing representative rows ---
$style1Rows = @(60,61,62,63,64,65,66,67,75,76,79,80,81,82,83)
$style2Rows = @(53,54,55,56,57,58,59,68,69,70,71,72,73,74,77,78)

$srcStyle1 = $ws.Range("A2:D2")
foreach ($r in $style1Rows) {
    $srcStyle1.Copy() | Out-Null
    $ws.Range("A$r" + ":D$r").PasteSpecial(-4122) | Out-Null
}

$srcStyle2 = $ws.Range("A13:D13")
foreach ($r in $style2Rows) {
    $srcStyle2.Copy() | Out-Null
    $ws.Range("A$r" + ":D$r").PasteSpecial(-4122) | Out-Null
}

# E column style (s=4) for rows 61 and 73, copy format from E16
$srcE = $ws.Range("E16")
foreach ($r in @(61,73)) {
    $srcE.Copy() | Out-Null
    $ws.Range("E$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Set cell values ---
$ws.Range("A53").Value = "Listado de los Servicio Contratados por el cliente"
$ws.Range("B53").Value = "dashboard/listado-publicaciones-contratadas"
$ws.Range("C53").Value = "Lista todos los servicios contratados por el cliente logueado."
$ws.Range("D53").Value = "OK."

$ws.Range("A54").Value = "Ver publicación contratada"
$ws.Range("B54").Value = "dashboard/listado-publicaciones-contratadas"
$ws.Range("C54").Value = "Click en el Link ver de una publicación."
$ws.Range("D54").Value = "OK. Redirecciona a dashboard/ver-publicacion-ofrecida/idPublicacionSeleccionada. Muestra todos los datos de la publicación."

$ws.Range("A55").Value = "Listado de los Servicio Contratados por el cliente"
$ws.Range("B55").Value = "dashboard/listado-publicaciones-contratadas"
$ws.Range("C55").Value = "Lista todos los servicios contratados por el cliente logueado."
$ws.Range("D55").Value = "OK."

$ws.Range("A56").Value = "Listado Comentarios Pendientes"
$ws.Range("B56").Value = "dashboard/listado-publicaciones-contratadas"
$ws.Range("C56").Value = "Lista todos los comentarios pendientes de los servicios contratados."
$ws.Range("D56").Value = "OK."

$ws.Range("A57").Value = "Comentar y Puntuar Cliente Contratado"
$ws.Range("B57").Value = "dashboard/listado-publicaciones-contratadas"
$ws.Range("C57").Value = "Click en el Link comentar. No se ingresa comentario."
$ws.Range("D57").Value = "Redirecciona a dashboard/ver-publicacion-ofrecida/idPublicacionSeleccionada. Muestra ventana modal para completar comentario y puntuación. Alert correspondiente indicando que se debe ingresar un comentario."

$ws.Range("A58").Value = "Comentar y Puntuar Cliente Contratado"
$ws.Range("B58").Value = "dashboard/listado-publicaciones-contratadas"
$ws.Range("C58").Value = "Click en el Link comentar. Se ingresa comentario de una letra."
$ws.Range("D58").Value = "Redirecciona a dashboard/ver-publicacion-ofrecida/idPublicacionSeleccionada. Muestra ventana modal para completar comentario y puntuación. Alert correspondiente indicando que se debe ingresar un comentario con al menos dos letras."

$ws.Range("A59").Value = "Comentar y Puntuar Cliente Contratado"
$ws.Range("B59").Value = "dashboard/listado-publicaciones-contratadas"
$ws.Range("C59").Value = "Click en el Link comentar. Comentario con dos o mas letras."
$ws.Range("D59").Value = "OK."

$ws.Range("A60").Value = "Datos del cliente"
$ws.Range("B60").Value = "dashboard/ver-perfil-usuario/idCliente"
$ws.Range("C60").Value = "Listado con los datos personales del cliente"
$ws.Range("D60").Value = "OK."

$ws.Range("A61").Value = "Puntaje por servicio del cliente"
$ws.Range("B61").Value = "dashboard/ver-perfil-usuario/idCliente"
$ws.Range("C61").Value = "Listado con los puntajes del cliente por servicio"
$ws.Range("D61").Value = "NO IMPLEMENTADO"
$ws.Range("E61").Value = "NO IMPLEMENTADO"

$ws.Range("A62").Value = "Comentario y puntuacion de los servicios ofrecidos por servicio del cliente"
$ws.Range("B62").Value = "dashboard/ver-perfil-usuario/idCliente"
$ws.Range("C62").Value = "Listado de todos los comentarios y puntuacion de los servicios ofrecidos por el cliente por servicio."
$ws.Range("D62").Value = "Alert correspondiente indicando que se debe ingresar un comentario."

$ws.Range("A63").Value = "Responder comentario de los servicios ofrecidos"
$ws.Range("B63").Value = "dashboard/ver-perfil-usuario/idCliente"
$ws.Range("C63").Value = "Click en link responder. No se ingresa comentario."
$ws.Range("D63").Value = "Alert correspondiente indicando que se debe ingresar un comentario."

$ws.Range("A64").Value = "Responder comentario de los servicios ofrecidos"
$ws.Range("B64").Value = "dashboard/ver-perfil-usuario/idCliente"
$ws.Range("C64").Value = "Click en link responder. Comentario correcto."
$ws.Range("D64").Value = "OK."

$ws.Range("A65").Value = "Comentario y puntuacion de las solicitudes realizadas por servicio del cliente"
$ws.Range("B65").Value = "dashboard/ver-perfil-usuario/idCliente"
$ws.Range("C65").Value = "Listado de todos los comentarios y puntuacion de las solicitudes realizadas por el cliente por servicio."
$ws.Range("D65").Value = "OK."

$ws.Range("A66").Value = "Responder comentario de las solicitudes realizadas"
$ws.Range("B66").Value = "dashboard/ver-perfil-usuario/idCliente"
$ws.Range("C66").Value = "Click en link responder. No se ingresa comentario."
$ws.Range("D66").Value = "Alert correspondiente indicando que se debe ingresar un comentario."

$ws.Range("A67").Value = "Responder comentario de las solicitudes realizadas"
$ws.Range("B67").Value = "dashboard/ver-perfil-usuario/idCliente"
$ws.Range("C67").Value = "Click en link responder. Comentario correcto."
$ws.Range("D67").Value = "OK."

$ws.Range("A68").Value = "Editar datos cliente"
$ws.Range("B68").Value = "dashboard/perfil-usuario"
$ws.Range("C68").Value = "Nombre, Apellido, teléfono, dirección, select barrio vacios."
$ws.Range("D68").Value = "Alert correspondiente indicando que debe completar cada dato."

$ws.Range("A69").Value = "Editar datos cliente"
$ws.Range("B69").Value = "dashboard/perfil-usuario"
$ws.Range("C69").Value = "Nombre y Apellido con menos de 2 caractéres."
$ws.Range("D69").Value = "Alert correspondiente indicando que debe ingresar un nombre y un apellido con al menos de dos caractéres."

$ws.Range("A70").Value = "Editar datos cliente"
$ws.Range("B70").Value = "dashboard/perfil-usuario"
$ws.Range("C70").Value = "Teléfono con menos de 6 caractéres."
$ws.Range("D70").Value = "Alert correspondiente indicando que debe ingresar un teléfono con al menos de 6 caratéres."

$ws.Range("A71").Value = "Editar datos cliente"
$ws.Range("B71").Value = "dashboard/perfil-usuario"
$ws.Range("C71").Value = "Dirección con menos de 4 caractéres."
$ws.Range("D71").Value = "Alert correspondiente indicando que debe ingresar una dirección con al menos de 4 caractéres."

$ws.Range("A72").Value = "Editar datos cliente"
$ws.Range("B72").Value = "dashboard/perfil-usuario"
$ws.Range("C72").Value = "Datos correctos"
$ws.Range("D72").Value = "OK. Alert indicando que los cambios se realizaron con éxito."

$ws.Range("A73").Value = "Editar datos cliente"
$ws.Range("B73").Value = "dashboard/perfil-usuario"
$ws.Range("C73").Value = "Click en selección de imagen."
$ws.Range("D73").Value = "NO IMPLEMENTADO"
$ws.Range("E73").Value = "NO IMPLEMENTADO"

$ws.Range("A74").Value = "Navegación"
$ws.Range("B74").Value = "dashboard/perfil-usuario"
$ws.Range("C74").Value = "Click en el Link cambiar contraseña."
$ws.Range("D74").Value = "OK. Redirecciona a dashboard/cambiar-contrasena-usuario."

$ws.Range("A75").Value = "Listado de todas las publicaciones de un servicio determinado"
$ws.Range("B75").Value = "dashboard/listado-servicios-ofrecidos/idServicio"
$ws.Range("C75").Value = "Listado de todas las publicaciones del servicio seleccionado."
$ws.Range("D75").Value = "OK. Muestra los datos de cada publicación y los datos del cliente dueño de la misma."

$ws.Range("A76").Value = "Navegación"
$ws.Range("B76").Value = "dashboard/listado-servicios-ofrecidos/idServicio"
$ws.Range("C76").Value = "Click en una publicación."
$ws.Range("D76").Value = "OK. Redirecciona a dashboard/ver-publicacion-ofrecida/idPublicacionSeleccionada. Muestra todos los datos de la publicación."

$ws.Range("A77").Value = "Listado de todas las solicitudes de un servicio determinado"
$ws.Range("B77").Value = "dashboard/listado-solicitudes-ofrecidas/idServicio"
$ws.Range("C77").Value = "Listado de todas las solicitudes del servicio seleccionado."
$ws.Range("D77").Value = "OK. Muestra los datos de cada publicación y los datos del cliente dueño de la misma."

$ws.Range("A78").Value = "Navegación"
$ws.Range("B78").Value = "dashboard/listado-solicitudes-ofrecidas/idServicio"
$ws.Range("C78").Value = "Click en una publicación."
$ws.Range("D78").Value = "OK. Redirecciona a dashboard/ver-publicacion-solicitada/idPublicacionSeleccionada. Muestra todos los datos de la publicación."

$ws.Range("A79").Value = "Cambiar contraseña"
$ws.Range("B79").Value = "dashboard/cambiar-contrasena-usuario"
$ws.Range("C79").Value = "Contraseña, nueva contraseña y repetir contraseña vacios."
$ws.Range("D79").Value = "Alert correspondientes indicando que debe ingresar la contraseña y la nueva contraseña."

$ws.Range("A80").Value = "Cambiar contraseña"
$ws.Range("B80").Value = "dashboard/cambiar-contrasena-usuario"
$ws.Range("C80").Value = "Contraseña anterior incorrecta."
$ws.Range("D80").Value = "Alert correspondientes indicando que hay un error en la contraseña anterior."

$ws.Range("A81").Value = "Cambiar contraseña"
$ws.Range("B81").Value = "dashboard/cambiar-contrasena-usuario"
$ws.Range("C81").Value = "Contraseña nueva con menos de 8 caracteres."
$ws.Range("D81").Value = "Alert correspondientes indicando que debe ingresar una nueva contraseña con al menos 8 caracteres."

$ws.Range("A82").Value = "Cambiar contraseña"
$ws.Range("B82").Value = "dashboard/cambiar-contrasena-usuario"
$ws.Range("C82").Value = "Contraseña nueva destinta de la confirmación."
$ws.Range("D82").Value = "Alert correspondientes indicando que la contraseña nueva y la confirmación no coinciden."

$ws.Range("A83").Value = "Cambiar contraseña"
$ws.Range("B83").Value = "dashboard/cambiar-contrasena-usuario"
$ws.Range("C83").Value = "Datos correctos."
$ws.Range("D83").Value = "OK."

# --- Column width adjustments ---
# Column A gets its own (new) width; column B keeps the old shared width but
# now lives in its own <col> entry since A and B no longer match.
$ws.Columns.Item(1).ColumnWidth = 67.65
# Column C: 70.109375 -> 99
$ws.Columns.Item(3).ColumnWidth = 98.1666666667
# Column D: 136.109375 -> 206
$ws.Columns.Item(4).ColumnWidth = 205.1666666667

# --- Update the active selection to the new "next empty row" ---
$ws.Range("A84").Select() | Out-Null
